$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = "足側板内側"
$ws.Range("G3").Value = "A2017"
$ws.Range("I3").Value = 1
$ws.Range("N3").Value = "無"
$ws.Range("O3").Value = "-"
$ws.Range("P3").Value = "-"
$ws.Range("B4").Value = "RR,RL"
$ws.Range("C4").Value = "R-P002"
$ws.Range("F4").Value = "足機構固定板"
$ws.Range("G4").Value = "AR-M2"
$ws.Range("I4").Value = 2
$ws.Range("N4").Value = "有"
$ws.Range("O4").Value = "アルマイト"
$ws.Range("P4").Value = "黒"
$ws.Range("B5").Value = "RR,RL"
$ws.Range("C5").Value = "R-P003"
$ws.Range("F5").Value = "足機構上面カバー"
$ws.Range("G5").Value = "AR-M2"
$ws.Range("B6").Value = "RR,RL"
$ws.Range("C6").Value = "R-P004"
$ws.Range("F6").Value = "足機構主柱"
$ws.Range("G6").Value = "A2017"
$ws.Range("H6").Value = 4
$ws.Range("I6").Value = 4
$ws.Range("C7").Value = "RR-P001"
$ws.Range("F7").Value = "足側板右外側"
$ws.Range("G7").Value = "A2017"
$ws.Range("C8").Value = "RR-P002"
$ws.Range("F8").Value = "右足機構側面充填樹脂"
$ws.Range("G8").Value = "ｴﾎﾟｷｼ樹脂(充填剤なし)"
$ws.Range("H8").Value = 1
$ws.Range("I8").Value = 1
$ws.Range("D11").Value = $null
$ws.Range("E11").Value = "MISUMI"
$ws.Range("F11").Value = "CB2.6-5"
$ws.Range("G11").Value = "SCM435"
$ws.Range("H11").Value = 4
$ws.Range("K11").Value = "六角穴付ﾎﾞﾙﾄ"
$ws.Range("L11").Value = "1日目"
$ws.Range("M11").Value = 642
$ws.Range("O11").Value = $null
$ws.Range("P11").Value = $null
$ws.Range("R11").Value = $null
$ws.Range("S11").Value = $null
$ws.Range("D12").Value = "Amazon"
$ws.Range("E12").Value = "マブチ"
$ws.Range("F12").Value = "RS380-PH"
$ws.Range("G12").Value = "材料 <指定なし>"
$ws.Range("H12").Value = 2
$ws.Range("J12").Value = "馬渕モータ RS380PH"
$ws.Range("K12").Value = "電動機"
$ws.Range("L12").Value = $null
$ws.Range("M12").Value = $null
$ws.Range("D13").Value = "KHK"
$ws.Range("E13").Value = "KHK"
$ws.Range("F13").Value = "SS0.5-10"
$ws.Range("G13").Value = "S45C"
$ws.Range("H13").Value = 2
$ws.Range("K13").Value = "平歯車"
$ws.Range("L13").Value = $null
$ws.Range("M13").Value = 150
$ws.Range("O13").Value = "-"
$ws.Range("P13").Value = "-"
$ws.Range("R13").Value = "-"
$ws.Range("S13").Value = "-"
$ws.Range("F14").Value = "CLJW6-8-40.0"
$ws.Range("G14").Value = "MCナイロン"
$ws.Range("H14").Value = 4
$ws.Range("K14").Value = "樹脂ｶﾗｰ"
$ws.Range("L14").Value = "2日目"
$ws.Range("M14").Value = 300
$ws.Range("F15").Value = "CLJW6-8-38.0"
$ws.Range("G15").Value = "MCナイロン"
$ws.Range("H15").Value = 4
$ws.Range("K15").Value = "樹脂ｶﾗｰ"
$ws.Range("L15").Value = "2日目"
$ws.Range("M15").Value = 300
$ws.Range("F16").Value = "CBSTBR3-6"
$ws.Range("H16").Value = 28
$ws.Range("K16").Value = "超極低頭ﾎﾞﾙﾄ"
$ws.Range("M16").Value = 280
$ws.Range("D17").Value = $null
$ws.Range("E17").Value = "MISUMI"
$ws.Range("F17").Value = "MPFZ10-8"
$ws.Range("G17").Value = "黄銅"
$ws.Range("J17").Value = $null
$ws.Range("K17").Value = "無給油ﾌﾞｯｼｭ"
$ws.Range("L17").Value = "1日目"
$ws.Range("M17").Value = 500
